# Automatic update from scheduled task
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the timestamp value already stored in A4 (tiny precision fix)
$ws.Range("A4").Value = 45864.3336262963

# Append the new reading captured by the scheduled task as row 5
$ws.Range("A5").Value = 45864.37523967567
$ws.Range("B5").Value = 2025
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 14
$ws.Range("E5").Value = 87.8
$ws.Range("F5").Value = 108.98
$ws.Range("G5").Value = 11.05
$ws.Range("H5").Value = "ESE"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "09:00:20"

# Match formatting of the other date cells in column A
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat
